$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: D1 becomes "median_percent" (was "median_1"); E1 ("median_2") will be removed below
$ws.Range("D1").Value = "median_percent"

# Update C and D values for each data row (row 2 - row 13)
$ws.Range("C2").Value = 0.01275
$ws.Range("D2").Value = 0.007634730538922156

$ws.Range("C3").Value = 0.0255
$ws.Range("D3").Value = 0.01526946107784431

$ws.Range("C4").Value = 0.138
$ws.Range("D4").Value = 0.08263473053892217

$ws.Range("C5").Value = 0.00035
$ws.Range("D5").Value = 0.0002095808383233533

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("C7").Value = 0.011
$ws.Range("D7").Value = 0.00658682634730539

$ws.Range("C8").Value = 0.002
$ws.Range("D8").Value = 0.001197604790419162

$ws.Range("C9").Value = 0.003
$ws.Range("D9").Value = 0.001796407185628743

$ws.Range("C10").Value = 0.091
$ws.Range("D10").Value = 0.05449101796407185

$ws.Range("C11").Value = 0.077
$ws.Range("D11").Value = 905.8823529411765

$ws.Range("C12").Value = 0.08995
$ws.Range("D12").Value = 1058.235294117647

$ws.Range("C13").Value = 0.33185
$ws.Range("D13").Value = 3904.117647058823

# Remove column E entirely (median_2 column), shrinking the used range to A1:D13
$ws.Range("E1:E13").Delete()
